# Applies the edits described by the commit "aanvulling op user permissions"
# to the Portfolio LU1 document.
#
# Strategy: use Find/Replace (wdFindContinue / wdReplaceOne) scoped over the
# whole document content for each textual change. Every search string below
# was checked to be unique in the document, so this is equivalent to editing
# each specific run/paragraph by hand.

$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $d.Content.Find.Execute($find, $false, $false, $false, $false, $false, `
                             $true, 1, $false, $replace, 2) | Out-Null
}

# --- 4. Acceptatiecriteria / Definition of Done bullet list -----------------

Replace-Text "het ERD de onderneming volledig en juist in kaart brengt." `
             "het ERD de onderneming volledig en juist in kaart gebracht zijn."

Replace-Text "de omzetting van ERD naar het Relationeel Model goed in elkaar overgaan. " `
             "de ERD goed omgezet is naar het Relationeel Model."

Replace-Text "de omzetting van het Relationeel Model en de Data Definition Language goed in elkaar overgaan." `
             "het Relationeel Model goed is omgezet naar de Data Definition Language."

Replace-Text "de DCL users db_admin, Medewerker, Bedrijf en Particulier aanmaakt." `
             "met DCL users db_admin, Medewerker, Bedrijf en Particulier aangemaakt zijn."

Replace-Text "de DCL wordt gecombineerd met Views die databasebeperkingen instellen voor de user Bedrijf en Particulier zodat zij geen data kunnen zien/bewerken van anderen." `
             "met DCL wordt gecombineerd met Views die databasebeperkingen instellen voor de user Bedrijf en Particulier zodat zij geen data kunnen zien/bewerken van anderen."

Replace-Text "de DCL rollen maakt voor meerdere users binnen de groep db_admin en Medewerker." `
             "met DCL rollen gemaakt zijn voor meerdere users binnen de groep db_admin en Medewerker."

Replace-Text "de DML worden uitgewerkt via vooraf bepaalde user stories." `
             "met DML de vooraf user stories uitgewerkt zijn."

Replace-Text "in de DML enkele stored procedures staan waarmee complexere berekeningen van de parkeerkosten gemaakt kunnen worden." `
             "met DML enkele Stored Procedures gemaakt zijn waarmee complexere berekeningen van de parkeerkosten gemaakt kunnen worden."

Replace-Text "in de DML events staan die blabla" `
             "met DML Events gemaakt zijn die blabla"

Replace-Text "in de DML triggers staan die blabla " `
             "met DML Triggers gemaakt zijn die blabla "

# --- Users section: add extra clarifying sentences on user permissions -----

Replace-Text "Oprichters en developers van PinPointParking" `
             "Oprichters en developers van PinPointParking die toegang hebben tot alle data en alle bewerkingen kunnen doen."

Replace-Text "Medewerkers van PinPointParking die administratieve taken hebben" `
             "Medewerkers van PinPointParking die administratieve taken hebben en alleen data kunnen zien en verwijderen vanwege de AVG."

Replace-Text "Zakelijke klant die medewerkers gebruik laten maken van PinPointParking, alle facturen gaan naar de financiële afdeling van het bedrijf." `
             "Zakelijke klant die medewerkers gebruik laat maken van PinPointParking, alle facturen gaan naar de financiële afdeling van het bedrijf. Zij kunnen alleen bij de data die samenhangt met hun eigen medewerkers. Deze kunnen ze zien, maken, aanpassen, verwijderen."

Replace-Text "Particuliere klanten die zelf hun facturen betalen." `
             "Particuliere klanten die zelf hun facturen betalen. Zij kunnen alleen bij de data die samenhangt met hun eigen account. Deze kunnen zij zien, maken, aanpassen, verwijderen."

# --- Move the _GoBack bookmark from "4. Acceptatiecriteria" heading --------
# --- down to the blank paragraph right after the Users list -----------------

if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# Find the specific blank paragraph that directly follows the
# "Particuliere klanten ..." paragraph (and precedes "User Stories medewerker"),
# and re-create the bookmark there.
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    $txt = $para.Range.Text
    if ($txt -like "*Particuliere klanten die zelf hun facturen betalen*") {
        $blankPara = $d.Paragraphs.Item($i + 1)
        $d.Bookmarks.Add("_GoBack", $blankPara.Range)
        break
    }
}
